$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item(2)

# Insert a new column before column A to hold the Year values
$ws.Columns.Item(1).Insert()

# Header row: Year column formula
$ws.Range("A1").Formula = "='Zambia Workbook'!A11"

# Data rows: Year values pulled from 'Zambia Workbook' column A (rows 12-32)
for ($r = 2; $r -le 22; $r++) {
    $srcRow = $r + 10
    $ws.Range("A$r").Formula = "='Zambia Workbook'!A$srcRow"
}

# New trailing columns H (iU / NA) and I (Country / ZM)
# Write values in the same order the shared-string table should end up in:
# iU, NA, Country, ZM
$ws.Range("H1").Value = "iU"
$ws.Range("H2").Value = "NA"
$ws.Range("I1").Value = "Country"
$ws.Range("I2").Value = "ZM"
for ($r = 3; $r -le 22; $r++) {
    $ws.Range("H$r").Value = "NA"
    $ws.Range("I$r").Value = "ZM"
}

Write-Host "Done"
